# Applies the cell value updates for the refreshed crypto symbol list
# (GitHub Actions data-refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe via .Formula so numeric-looking strings
# (e.g. "0.1680", "0.94%") are stored as literal text, matching the
# original inlineStr cells, instead of being parsed into numbers/percentages.
$ws.Range("D2").Formula = "'294.24"
$ws.Range("E2").Formula = "'0.94%"
$ws.Range("D3").Formula = "'31.08"
$ws.Range("E3").Formula = "'0.53%"
$ws.Range("D4").Formula = "'4.932"
$ws.Range("E4").Formula = "'1.19%"
$ws.Range("D5").Formula = "'0.07351"
$ws.Range("E5").Formula = "'1.88%"
$ws.Range("D6").Formula = "'2.306"
$ws.Range("E6").Formula = "'33.38%"
$ws.Range("D7").Formula = "'7.729"
$ws.Range("E7").Formula = "'0.79%"
$ws.Range("D8").Formula = "'3.741"
$ws.Range("E8").Formula = "'-0.57%"
$ws.Range("D9").Formula = "'0.9061"
$ws.Range("E9").Formula = "'0.33%"
$ws.Range("D10").Formula = "'0.1680"
$ws.Range("E10").Formula = "'1.71%"
$ws.Range("D11").Formula = "'0.07994"
$ws.Range("E11").Formula = "'5.43%"
$ws.Range("D12").Formula = "'0.08137"
$ws.Range("E12").Formula = "'1.66%"
$ws.Range("D13").Formula = "'0.03098"
$ws.Range("E13").Formula = "'1.75%"
$ws.Range("D14").Formula = "'0.1009"
$ws.Range("E14").Formula = "'0.92%"
$ws.Range("D15").Formula = "'0.001515"
$ws.Range("E15").Formula = "'1.33%"
$ws.Range("E16").Formula = "'2.72%"
$ws.Range("D17").Formula = "'3.478"
$ws.Range("E17").Formula = "'0.50%"
$ws.Range("E18").Formula = "'-1.49%"
$ws.Range("E19").Formula = "'1.14%"
$ws.Range("D20").Formula = "'0.1303"
$ws.Range("E20").Formula = "'-0.09%"
$ws.Range("D21").Formula = "'3.970"
$ws.Range("E21").Formula = "'-9.75%"
$ws.Range("D22").Formula = "'0.2097"
$ws.Range("E22").Formula = "'4.68%"
$ws.Range("D23").Formula = "'0.04528"
$ws.Range("E23").Formula = "'0.90%"
$ws.Range("D24").Formula = "'0.001212"
$ws.Range("E24").Formula = "'-0.57%"
$ws.Range("D25").Formula = "'0.004656"
$ws.Range("E25").Formula = "'15.72%"
$ws.Range("D26").Formula = "'0.0001300"
$ws.Range("E26").Formula = "'3.46%"
$ws.Range("D27").Formula = "'0.0003389"
$ws.Range("D39").Formula = "'0.01614"
$ws.Range("E39").Formula = "'-2.64%"
$ws.Range("E40").Formula = "'2.47%"
$ws.Range("D41").Formula = "'0.007361"
$ws.Range("E41").Formula = "'-1.32%"
$ws.Range("B42").Formula = "'Dexo"
$ws.Range("C42").Formula = "'https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Formula = "'0.008784"
$ws.Range("E42").Formula = "'--%"
$ws.Range("B43").Formula = "'BKEXToken"
$ws.Range("C43").Formula = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Formula = "'0.1332"
$ws.Range("E43").Formula = "'1.25%"
$ws.Range("D44").Formula = "'0.002001"
$ws.Range("E44").Formula = "'-2.30%"
$ws.Range("D45").Formula = "'0.009512"
$ws.Range("E45").Formula = "'-6.96%"
$ws.Range("D46").Formula = "'0.00005919"
$ws.Range("E46").Formula = "'3.49%"
$ws.Range("D47").Formula = "'0.00000000750"
$ws.Range("E47").Formula = "'-0.51%"
$ws.Range("E48").Formula = "'3.02%"
$ws.Range("D49").Formula = "'0.002894"
$ws.Range("E49").Formula = "'-4.02%"
$ws.Range("D50").Formula = "'0.00002100"
$ws.Range("E50").Formula = "'-0.51%"
$ws.Range("D51").Formula = "'0.0002000"
$ws.Range("E51").Formula = "'-0.51%"
